# Swap the record data between row 2 <-> row 3 and row 8 <-> row 9.
# Only the columns that actually differ between each pair of rows are
# touched (A:B, E:H, Q:R, Z, AB for the 2/3 pair; A:B, D:H, Q:R, Z, AB for
# the 8/9 pair) so that untouched columns (dates in Y/AA, location text,
# etc.) are left completely alone and avoid any implicit type coercion
# (e.g. Excel turning a "YYYY-MM-DD" string into a date serial on write).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Range {
    param($ws, [string]$addr1, [string]$addr2)
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# --- Row 2 <-> Row 3 ---
Swap-Range $ws "A2:B2" "A3:B3"
Swap-Range $ws "E2:H2" "E3:H3"
Swap-Range $ws "Q2:R2" "Q3:R3"
Swap-Range $ws "Z2"    "Z3"
Swap-Range $ws "AB2"   "AB3"

# --- Row 8 <-> Row 9 ---
Swap-Range $ws "A8:B8" "A9:B9"
Swap-Range $ws "D8:H8" "D9:H9"
Swap-Range $ws "Q8:R8" "Q9:R9"
Swap-Range $ws "Z8"    "Z9"
Swap-Range $ws "AB8"   "AB9"
